# "Add files via upload" — the re-uploaded workbook gained two new cells
# on row 1 of the "Schedule" sheet: C1 = "dadaw" and E1 = "32423" (kept as
# text, not a number, matching the shared-string cell in the source file).
# Nothing else on the sheet changes; Excel automatically extends the used
# range / dimension to include the new row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "dadaw"

# Leading apostrophe forces Excel to store this numeric-looking value as
# text (shared string), exactly like the source workbook does.
$ws.Range("E1").Value = "'32423"
